$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new rows 20-22 with new reliable species
$ws.Range("A20").Value = "PECTJAC"
$ws.Range("A21").Value = "SCYOCAN"
$ws.Range("A22").Value = "RAJAAST"

# Row 12: species FLEXGLAB -> CHLAGLA (source column E12 stays "fishbase")
$ws.Range("A12").Value = "CHLAGLA"

# Update selection to match the saved view state
$ws.Range("A24").Select()
